# Updates cryptos list (prices / 1h volume %, plus a couple of row
# re-ranks where two coins swapped places) per the Dec 22 2023 data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "44.217.31"
$ws.Cells.Item(2, 5).Value = "  +1.26%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.276.41"
$ws.Cells.Item(3, 5).Value = "  +3.27%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.04%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "97.93"
$ws.Cells.Item(5, 5).Value = "  +15.07%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "275.30"
$ws.Cells.Item(6, 5).Value = "  +6.55%  "

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.629"
$ws.Cells.Item(7, 5).Value = "  +1.76%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  +0.03%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.643"
$ws.Cells.Item(9, 5).Value = "  +8.17%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "48.06"
$ws.Cells.Item(10, 5).Value = "  +7.21%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0945"
$ws.Cells.Item(11, 5).Value = "  +3.29%  "

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "8.20"
$ws.Cells.Item(12, 5).Value = "  +11.18%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.105"
$ws.Cells.Item(13, 5).Value = "  +1.09%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "15.75"
$ws.Cells.Item(14, 5).Value = "  +9.78%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "2.607.77"
$ws.Cells.Item(15, 5).Value = "  +2.94%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "0.843"
$ws.Cells.Item(16, 5).Value = "  +7.81%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "2.257.35"
$ws.Cells.Item(17, 5).Value = "  +2.15%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "44.207.09"
$ws.Cells.Item(18, 5).Value = "  +1.32%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "  +3.37%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "6.25"
$ws.Cells.Item(20, 5).Value = "  +5.65%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "71.11"
$ws.Cells.Item(21, 5).Value = "  +1.90%  "

# Row 22
$ws.Cells.Item(22, 2).Value = "ImmutableX"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "2.33"
$ws.Cells.Item(22, 5).Value = "  +0.26%  "

# Row 23
$ws.Cells.Item(23, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "10.27"
$ws.Cells.Item(23, 5).Value = "  +14.98%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "235.40"
$ws.Cells.Item(24, 5).Value = "  +1.94%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "1.00"
$ws.Cells.Item(25, 5).Value = "  -0.03%  "

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "11.58"
$ws.Cells.Item(26, 5).Value = "  +9.06%  "

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "2.54"
$ws.Cells.Item(27, 5).Value = "  +13.42%  "

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "39.68"
$ws.Cells.Item(28, 5).Value = "  +0.77%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "3.35"
$ws.Cells.Item(29, 5).Value = "  -5.84%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  +0.13%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "173.68"
$ws.Cells.Item(31, 5).Value = "  -0.03%  "

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.0925"
$ws.Cells.Item(32, 5).Value = "  +6.85%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "21.40"
$ws.Cells.Item(33, 5).Value = "  +4.91%  "

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "5.70"
$ws.Cells.Item(34, 5).Value = "  +6.94%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  +1.96%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.114"
$ws.Cells.Item(36, 5).Value = "  +3.74%  "

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.0358"
$ws.Cells.Item(37, 5).Value = "  -0.37%  "

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "4.43"
$ws.Cells.Item(38, 5).Value = "  -0.66%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "3.60"
$ws.Cells.Item(39, 5).Value = "  +26.42%  "

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.257"
$ws.Cells.Item(40, 5).Value = "  +29.82%  "

# Row 41
$ws.Cells.Item(41, 2).Value = "LidoDAOToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "2.21"
$ws.Cells.Item(41, 5).Value = "  +5.50%  "

# Row 42
$ws.Cells.Item(42, 2).Value = "Celestia"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "12.56"
$ws.Cells.Item(42, 5).Value = "  +0.64%  "

# Row 43
$ws.Cells.Item(43, 2).Value = "THORChain"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "5.52"
$ws.Cells.Item(43, 5).Value = "  +1.05%  "

# Row 44
$ws.Cells.Item(44, 2).Value = "MultiversX"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "62.96"
$ws.Cells.Item(44, 5).Value = "  -0.28%  "

# Row 45
$ws.Cells.Item(45, 5).Value = "  +5.43%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "8.60"
$ws.Cells.Item(46, 5).Value = "  +3.04%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "100.62"
$ws.Cells.Item(47, 5).Value = "  +0.27%  "

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "1.17"
$ws.Cells.Item(48, 5).Value = "  +6.08%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "  +1.61%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.437"
$ws.Cells.Item(50, 5).Value = "  -0.30%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "2.492.70"
$ws.Cells.Item(51, 5).Value = "  +3.10%  "
